# "Finalizado os testes da pesquisa"
#
# The author finished reviewing the results table: the two data columns
# (A = dataset name, B = accuracy) were widened to fit their contents
# (Format > AutoFit Column Width), and the sheet's scroll/selection state
# was reset to just past the last row of data (A43) instead of sitting on
# the previously-selected H41 with the view scrolled down to row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("dataset") and column B ("acc") are auto-fit to their longest
# entries. ColumnWidth is expressed in characters; the values below are the
# character widths whose on-disk (1/256-char-quantized) stored width lands
# closest to the real AutoFit result Excel computed for these columns
# (30.5703125 for A, 26.28515625 for B).
$ws.Columns.Item(1).ColumnWidth = 29.666666666666664
$ws.Columns.Item(2).ColumnWidth = 25.5

# Reset the selection/scroll position: select A43 (just below the last
# data row), which also clears the previous topLeftCell scroll anchor.
$ws.Range("A43").Select()
